$wb = $excel.ActiveWorkbook
$kobe = $wb.Worksheets.Item("kobe")

# --- Update existing daily figures (rows 57-60) ---------------------------
$kobe.Cells.Item(57, 6).Value = 52   # F57
$kobe.Cells.Item(58, 6).Value = 57   # F58

$kobe.Cells.Item(59, 4).Value = 21   # D59
$kobe.Cells.Item(59, 5).Value = 120  # E59
$kobe.Cells.Item(59, 6).Value = 57   # F59
$kobe.Cells.Item(59, 7).Value = 54   # G59

$kobe.Cells.Item(60, 2).Value = 47   # B60
$kobe.Cells.Item(60, 3).Value = 790  # C60
$kobe.Cells.Item(60, 4).Value = 7    # D60
$kobe.Cells.Item(60, 5).Value = 127  # E60
$kobe.Cells.Item(60, 9).Value = 1    # I60

# --- Insert a new row for the next day (row 61) ----------------------------
$kobe.Rows.Item(61).Insert()

$kobe.Cells.Item(61, 1).Value = 43934  # date (2020-04-13)
$kobe.Cells.Item(61, 2).Value = 20
$kobe.Cells.Item(61, 3).Value = 810
$kobe.Cells.Item(61, 4).Value = 1
$kobe.Cells.Item(61, 5).Value = 128
$kobe.Cells.Item(61, 6).Value = 56
$kobe.Cells.Item(61, 7).Value = 53
$kobe.Cells.Item(61, 8).Value = 3
$kobe.Cells.Item(61, 9).Value = 1
$kobe.Cells.Item(61, 10).Value = 26

# --- Update the print area to include the extra row ------------------------
$kobe.PageSetup.PrintArea = '$A$1:$J$64'

# --- Restore the view: kobe becomes the active/selected sheet --------------
$kobe.Activate()
$kobe.Application.ActiveWindow.SelectedSheets.Item(1)
$kobe.Range("B61").Select()
